# Refresh the crypto price/volume table with the latest scraped values.
# (GitHub Actions scheduled update.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    # Write the value as plain text, preserving representations like
    # trailing zeros ("575.10") or thousand-dot-separated numbers
    # ("2.462.63") that Excel would otherwise silently reinterpret as
    # numeric values, then restore the default (unformatted) cell
    # style so no stray formatting is left behind.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextCell "D2" "62.826.67"
Set-TextCell "E2" "  +0.12%  "
Set-TextCell "D3" "2.462.63"
Set-TextCell "E3" "  +0.70%  "
Set-TextCell "E4" "  +0.04%  "
Set-TextCell "D5" "575.10"
Set-TextCell "E5" "  -0.20%  "
Set-TextCell "D6" "147.10"
Set-TextCell "E6" "  +0.95%  "
Set-TextCell "E7" "  -0.03%  "
Set-TextCell "E8" "  -0.90%  "
Set-TextCell "D9" "2.462.46"
Set-TextCell "E9" "  +0.75%  "
Set-TextCell "E10" "  +0.68%  "
Set-TextCell "E11" "  -0.30%  "
Set-TextCell "E12" "  +0.06%  "
Set-TextCell "D13" "0.356"
Set-TextCell "E13" "  +0.86%  "
Set-TextCell "D14" "29.02"
Set-TextCell "E14" "  +2.29%  "
Set-TextCell "E15" "  -0.75%  "
Set-TextCell "D16" "2.910.42"
Set-TextCell "E16" "  +0.73%  "
Set-TextCell "D17" "62.737.60"
Set-TextCell "E17" "  +0.19%  "
Set-TextCell "D18" "2.466.77"
Set-TextCell "E18" "  +1.18%  "
Set-TextCell "E19" "  +0.15%  "
Set-TextCell "D20" "11.03"
Set-TextCell "E20" "  +0.11%  "
Set-TextCell "D21" "326.37"
Set-TextCell "E21" "  -1.16%  "
Set-TextCell "E22" "  +8.58%  "
Set-TextCell "E23" "  -0.08%  "
Set-TextCell "E24" "  -0.01%  "
Set-TextCell "D25" "10.07"
Set-TextCell "E25" "  +17.62%  "
Set-TextCell "D26" "65.58"
Set-TextCell "E26" "  -1.08%  "
Set-TextCell "D27" "643.76"
Set-TextCell "E27" "  -0.65%  "
Set-TextCell "D28" "0.0₃0983"
Set-TextCell "E28" "  -0.73%  "
Set-TextCell "D29" "2.582.90"
Set-TextCell "E29" "  +0.67%  "
Set-TextCell "D30" "0.999"
Set-TextCell "E30" "  -15.14%  "
Set-TextCell "E31" "  -0.66%  "
Set-TextCell "D32" "7.98"
Set-TextCell "E32" "  -2.54%  "
Set-TextCell "E33" "  -1.35%  "
Set-TextCell "E34" "  -2.87%  "
Set-TextCell "E35" "  -0.05%  "
Set-TextCell "E36" "  +2.72%  "
Set-TextCell "D37" "4.75"
Set-TextCell "E37" "  -0.27%  "
Set-TextCell "B38" "Monero"
Set-TextCell "C38" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D38" "152.14"
Set-TextCell "E38" "  -0.84%  "
Set-TextCell "B39" "PolygonEcosystemToken"
Set-TextCell "C39" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextCell "D39" "0.369"
Set-TextCell "E39" "  -1.40%  "
Set-TextCell "B40" "dogwifhat"
Set-TextCell "C40" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D40" "2.80"
Set-TextCell "E40" "  +2.95%  "
Set-TextCell "D41" "18.70"
Set-TextCell "E41" "  -0.30%  "
Set-TextCell "E42" "  -2.18%  "
Set-TextCell "E44" "  -37.04%  "
Set-TextCell "E45" "  -0.03%  "
Set-TextCell "D46" "152.42"
Set-TextCell "E46" "  +4.99%  "
Set-TextCell "E47" "  +2.22%  "
Set-TextCell "E48" "  -1.29%  "
Set-TextCell "E49" "  -0.59%  "
Set-TextCell "D50" "0.606"
Set-TextCell "E50" "  +0.17%  "
Set-TextCell "E51" "  -0.96%  "
